$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.831
$ws.Range("D9").Value = -7.545
$ws.Range("E12").Value = 17.646
$ws.Range("D13").Value = -8.263000000000002
$ws.Range("E14").Value = 17.007
$ws.Range("D16").Value = -8.461000000000002
$ws.Range("D18").Value = -8.547999999999998
$ws.Range("E19").Value = 16.566
$ws.Range("D20").Value = -7.438
$ws.Range("D26").Value = -7.462000000000001
$ws.Range("E26").Value = 17.129
$ws.Range("D27").Value = -8.055999999999999
$ws.Range("E27").Value = 16.925
$ws.Range("D29").Value = -7.415999999999999
$ws.Range("E29").Value = 16.904
$ws.Range("D35").Value = -7.606
$ws.Range("D36").Value = -7.783999999999999
$ws.Range("E37").Value = 16.792
$ws.Range("E38").Value = 16.835
$ws.Range("D45").Value = -7.629
$ws.Range("E47").Value = 16.77
$ws.Range("E51").Value = 16.558
$ws.Range("E52").Value = 16.8
$ws.Range("D55").Value = -8.224
$ws.Range("E55").Value = 16.669
$ws.Range("D57").Value = -8.315000000000001
$ws.Range("D69").Value = -7.470999999999999
$ws.Range("E69").Value = 17.255
$ws.Range("E70").Value = 17.609
$ws.Range("D76").Value = -7.806999999999999
$ws.Range("E76").Value = 16.732
$ws.Range("D78").Value = -8.123999999999999
$ws.Range("E81").Value = 16.457
$ws.Range("D82").Value = -8.418000000000001
$ws.Range("D83").Value = -8.031000000000001
$ws.Range("E83").Value = 16.767
$ws.Range("D93").Value = -7.452
$ws.Range("E94").Value = 18.004
$ws.Range("D97").Value = -8.184000000000001
$ws.Range("E100").Value = 16.489
$ws.Range("E102").Value = 16.701
